$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2-10, column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 124
$ws1.Range("F3").Value = 260
$ws1.Range("F4").Value = 140
$ws1.Range("F5").Value = 1740
$ws1.Range("F6").Value = 1507
$ws1.Range("F7").Value = 276
$ws1.Range("F8").Value = 62
$ws1.Range("F9").Value = 470
$ws1.Range("F10").Value = 121

# Sheet "全部类型" (all types) - rows 2-11 (row 8 is a different event, unchanged), column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 124
$ws4.Range("F3").Value = 260
$ws4.Range("F4").Value = 140
$ws4.Range("F5").Value = 1740
$ws4.Range("F6").Value = 1507
$ws4.Range("F7").Value = 276
$ws4.Range("F9").Value = 62
$ws4.Range("F10").Value = 470
$ws4.Range("F11").Value = 121
